# Auto-generated Excel COM-interop edit script
# Updates cryptos list price/volume(1h) columns, including two pairs of
# reordered rows (37/38, 41/42, 49/50) per upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.343.84"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "3.111.87"
$ws.Range("E3").Value = "  -4.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.91"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.67"
$ws.Range("E6").Value = "  +4.81%  "

$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("D8").Value = "3.106.19"
$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +0.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.95"
$ws.Range("E11").Value = "  +2.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.96"
$ws.Range("E13").Value = "  +2.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000244"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("D15").Value = "3.620.96"
$ws.Range("E15").Value = "  -3.70%  "

$ws.Range("E16").Value = "  -1.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.25"
$ws.Range("E17").Value = "  +2.70%  "

$ws.Range("D18").Value = "63.890.33"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").Value = "3.104.87"
$ws.Range("E19").Value = "  -2.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.27"
$ws.Range("E20").Value = "  +1.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.92"
$ws.Range("E21").Value = "  +4.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.739"
$ws.Range("E22").Value = "  +1.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.61"
$ws.Range("E23").Value = "  +2.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.36"
$ws.Range("E24").Value = "  +3.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  +7.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.85"
$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("E28").Value = "  +7.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.50"
$ws.Range("E29").Value = "  +5.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.71"
$ws.Range("E30").Value = "  +1.23%  "

$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.21"
$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("E33").Value = "  +6.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.53"
$ws.Range("E34").Value = "  +1.99%  "

$ws.Range("D35").Value = "0.0₃0851"
$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("E36").Value = "  +1.50%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.40"
$ws.Range("E37").Value = "  +2.58%  "

$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.17"
$ws.Range("E38").Value = "  +2.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.28"
$ws.Range("E39").Value = "  -1.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.39"
$ws.Range("E40").Value = "  +6.50%  "

$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "457.59"
$ws.Range("E41").Value = "  +4.67%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "50.86"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.291"
$ws.Range("E43").Value = "  +3.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0371"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D45").Value = "2.853.78"
$ws.Range("E45").Value = "  -1.96%  "

$ws.Range("E46").Value = "  +1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.05"
$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.47"
$ws.Range("E48").Value = "  +2.37%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.35"
$ws.Range("E49").Value = "  +5.02%  "

$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.28"
$ws.Range("E51").Value = "  +4.28%  "

